# This script reproduces (on the live Word object model) the human edit that is
# described by the target diff:
#
#   Paragraph "... كه برابر با "آنچه نتيجه ميدهد" است. در زبان پارسي، ميتوان
#   از "نتيجه" استفاده کرد."
#       -> the phrase "آنچه نتيجه ميدهد" is retyped as "سرانجام", and the
#          trailing sentence ("در زبان پارسي، ...") is deleted outright.
#       -> because this is now the location of the most recent edit, Word's
#          "_GoBack" bookmark ends up sitting at the end of that paragraph.
#
#   Paragraph ""صحت و سقم" >> "..."
#       -> this paragraph used to hold the "_GoBack" bookmark (from a
#          previous edit session); since the bookmark moves away, the two
#          runs that used to be split by it collapse back into a single run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the two paragraphs we need to touch by their (stable) text content,
# rather than by a hard-coded paragraph index.
# ---------------------------------------------------------------------------
$paraWithMahsel = $null
$paraWithSehhat = $null
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $paragraphs.Item($i)
    $text = $candidate.Range.Text
    if ($text.Contains("ماحصل")) {
        $paraWithMahsel = $candidate
    }
    if ($text.Contains("صحت و سقم")) {
        $paraWithSehhat = $candidate
    }
}

# ---------------------------------------------------------------------------
# Step 1: inside the "ماحصل" paragraph, retype "آنچه نتیجه میدهد" as "سرانجام"
# ---------------------------------------------------------------------------
$paraRange = $d.Range($paraWithMahsel.Range.Start, $paraWithMahsel.Range.End)
$paraRange.Find.ClearFormatting()
$paraRange.Find.Execute("آنچه نتیجه میدهد", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$phraseStart = $paraRange.Start
$phraseEnd = $paraRange.End

$replaceRange = $d.Range($phraseStart, $phraseEnd)
$replaceRange.Text = "سرانجام"
$replacedEnd = $phraseStart + 7   # length of "سرانجام"

# Force the freshly-typed word to live in its own run (instead of being
# silently re-merged with the identical formatting of its neighbour) by
# toggling a character property on and back off.
$freezeLeft = $d.Range($phraseStart, $replacedEnd)
$freezeLeft.Bold = 1
$freezeLeft.Bold = 0

# ---------------------------------------------------------------------------
# Step 2: move the "_GoBack" bookmark from the "صحت و سقم" paragraph to the
# end of the "ماحصل" paragraph (right after the closing "... است.").
# ---------------------------------------------------------------------------
$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}

# '" است.' is exactly 6 characters long and is the text that must remain
# right before the bookmark.
$bookmarkPos = $replacedEnd + 6
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$bookmarks.Add("_GoBack", $bookmarkRange)

# ---------------------------------------------------------------------------
# Step 3: delete the now-orphaned trailing sentence that used to follow
# '" است.' (i.e. everything up to, but not including, the paragraph mark).
# ---------------------------------------------------------------------------
$paraEndNow = $paraWithMahsel.Range.End
$tailRange = $d.Range($bookmarkPos, $paraEndNow - 1)
if ($tailRange.Start -lt $tailRange.End) {
    $tailRange.Delete()
}

# Re-freeze the same boundary so the trailing '" است.' keeps living in its
# own run rather than merging back into the "سرانجام" run now that the text
# after it is gone.
$freezeLeft2 = $d.Range($phraseStart, $replacedEnd)
$freezeLeft2.Bold = 1
$freezeLeft2.Bold = 0

# ---------------------------------------------------------------------------
# Step 4: now that the bookmark is gone from the "صحت و سقم" paragraph, the
# two runs that used to straddle it would naturally be saved back as a
# single run. Force that normalization explicitly: briefly change the text
# then restore it, which collapses same-formatted adjacent runs.
# ---------------------------------------------------------------------------
$sehhatStart = $paraWithSehhat.Range.Start
$sehhatEnd = $paraWithSehhat.Range.End
$sehhatRange = $d.Range($sehhatStart, $sehhatEnd - 1)
$originalSehhatText = $sehhatRange.Text

$sehhatRange.Text = "placeholder-text-for-run-merge"
$sehhatEndNow = $paraWithSehhat.Range.End
$sehhatRangeFinal = $d.Range($sehhatStart, $sehhatEndNow - 1)
$sehhatRangeFinal.Text = $originalSehhatText
